# Saldo_guide.xlsx update
# - Rename worksheet tab to reflect the new extraction timestamp
# - Shift the "Dt. Referencia" (column G) date for every data row from
#   2024-06-17 (serial 45460) to 2024-06-18 (serial 45461)
# - Correct a few "Saldo Previsto" / "Vl. Projetado" / "Vl. Total" values
#   that were revised for rows 17, 101 and 103

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the new run timestamp embedded in its name
$ws.Name = "IClientBalance-20240618-094409-"

# Shift every row's reference date (column G) by one day: 45460 -> 45461
$lastRow = 257
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45461
}

# Row 17: Saldo Previsto (D) and Vl. Total (H) corrected
$ws.Cells.Item(17, 4).Value = 2925.59
$ws.Cells.Item(17, 8).Value = 2925.59

# Row 101: Saldo Previsto (D) and Vl. Total (H) corrected
$ws.Cells.Item(101, 4).Value = 386.67
$ws.Cells.Item(101, 8).Value = 386.67

# Row 103: Saldo Previsto (D) corrected and Vl. Projetado (E) zeroed out
# (Vl. Total (H) for this row was already 3015.66 and remains unchanged)
$ws.Cells.Item(103, 4).Value = 3015.66
$ws.Cells.Item(103, 5).Value = 0
